$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "06/29/2025"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 0.0004638799999999998
$ws.Cells.Item($row, 3).Value = 107786.4965077175
$ws.Cells.Item($row, 4).Value = 50
